$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows involved in the rotation (data rows 195, 196, 197 on the worksheet)
$rows = @(195, 196, 197)

# Columns that participate in the rotation: B, and E through AD (A, C, D stay fixed)
$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# Capture the current ("before") values for each row/column combination
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{}
    foreach ($c in $cols) {
        $data[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Rotate: new row195 <- old row196, new row196 <- old row197, new row197 <- old row195
$mapping = @{ 195 = 196; 196 = 197; 197 = 195 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $data[$src][$c]
    }
}
